# FINFLUX-3612 Cartias specific scenarios
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 506.05
$wsSummary.Range("E2").Value = 4493.95
$wsSummary.Range("F3").Value = 3.09
$wsSummary.Range("A5").Value = 0
$wsSummary.Range("B5").Value = 0
$wsSummary.Range("A1:A1048576").Select()

# ---- Original Schedule sheet ----
$wsOriginal = $wb.Worksheets.Item("Original Schedule")
$wsOriginal.Range("F4").Value = 0
$wsOriginal.Range("G4").Value = 46.03

# ---- Repayment schedule sheet ----
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("J4").Value = 0
$wsRepay.Range("K4").Value = 42.94
$wsRepay.Range("Q4").Value = 42.94
$wsRepay.Range("J19").Select()

# ---- Transactions sheet ----
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 661
$wsTrans.Range("J2").Value = 4533.79
$wsTrans.Range("A3").Value = 660
$wsTrans.Range("A4").Value = 659
$wsTrans.Range("E4").Value = 43.67
$wsTrans.Range("J4").Value = 4488.05
$wsTrans.Range("A5").Value = 658
$wsTrans.Range("E5").Value = 43.67
$wsTrans.Range("A6").Value = 657
$wsTrans.Range("A7").Value = 652
$wsTrans.Range("A8").Value = 651
$wsTrans.Range("A9").Value = 650
$wsTrans.Range("H5").Select()

# ---- ChargesTab sheet ----
$wsCharges = $wb.Worksheets.Item("ChargesTab")
$wsCharges.Rows.Item(2).Delete()
$wsCharges.Rows.Item(2).Delete()
$wsCharges.Range("F9").Select()
